# Update "想去人数" (F column) values across sheets, matching the
# regenerated data snapshot (gh-pages output at 3967b19).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F12").Value = 711
$wsExhibit.Range("F13").Value = 1203
$wsExhibit.Range("F19").Value = 348
$wsExhibit.Range("F24").Value = 709
$wsExhibit.Range("F26").Value = 133
$wsExhibit.Range("F27").Value = 343

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F8").Value = 121
$wsShow.Range("F9").Value = 52
$wsShow.Range("F10").Value = 460

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F16").Value = 711
$wsAll.Range("F17").Value = 1203
$wsAll.Range("F23").Value = 348
$wsAll.Range("F29").Value = 121
$wsAll.Range("F30").Value = 121
$wsAll.Range("F33").Value = 709
$wsAll.Range("F35").Value = 52
$wsAll.Range("F36").Value = 133
$wsAll.Range("F37").Value = 460
$wsAll.Range("F39").Value = 343
